$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update scripts with new TPM values (recomputed ligand/receptor/edge stats)

# Row 2 (ECs -> ECs)
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.04155
$ws.Range("H2").Value = 0.12465
$ws.Range("I2").Value = 0.0001466168179836329
$ws.Range("J2").Value = 0.0001466168179836329
$ws.Range("M2").Value = 123.2806423333333
$ws.Range("N2").Value = 369.841927
$ws.Range("O2").Value = 0.6241574062367528
$ws.Range("P2").Value = 0.6241574062367526
$ws.Range("Q2").Value = 5.12231068895
$ws.Range("R2").Value = 46.10079620055
$ws.Range("S2").Value = 0.00009151197282335037
$ws.Range("T2").Value = 0.00009151197282335036

# Row 3 (ECs -> FAPs)
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.04155
$ws.Range("H3").Value = 0.12465
$ws.Range("I3").Value = 0.0001466168179836329
$ws.Range("J3").Value = 0.0001466168179836329
$ws.Range("O3").Value = 0.2392728888301323
$ws.Range("P3").Value = 0.2392728888301322
$ws.Range("Q3").Value = 1.963655423749999
$ws.Range("R3").Value = 17.67289881375
$ws.Range("S3").Value = 0.00003508142959002552
$ws.Range("T3").Value = 0.00003508142959002552

# Row 4 (ECs -> MuSCs)
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.04155
$ws.Range("H4").Value = 0.12465
$ws.Range("I4").Value = 0.0001466168179836329
$ws.Range("J4").Value = 0.0001466168179836329
$ws.Range("O4").Value = 0.136569704933115
$ws.Range("P4").Value = 0.136569704933115
$ws.Range("Q4").Value = 1.1207949347
$ws.Range("R4").Value = 10.0871544123
$ws.Range("S4").Value = 0.00002002341557025698
$ws.Range("T4").Value = 0.00002002341557025697

# Row 5 (FAPs -> ECs)
$ws.Range("I5").Value = 0.9992428949822291
$ws.Range("J5").Value = 0.9992428949822291
$ws.Range("M5").Value = 123.2806423333333
$ws.Range("N5").Value = 369.841927
$ws.Range("O5").Value = 0.6241574062367528
$ws.Range("P5").Value = 0.6241574062367526
$ws.Range("Q5").Value = 34910.26904155153
$ws.Range("R5").Value = 314192.4213739638
$ws.Range("S5").Value = 0.623684853532612
$ws.Range("T5").Value = 0.6236848535326119

# Row 6 (FAPs -> FAPs)
$ws.Range("I6").Value = 0.9992428949822291
$ws.Range("J6").Value = 0.9992428949822291
$ws.Range("O6").Value = 0.2392728888301323
$ws.Range("P6").Value = 0.2392728888301322
$ws.Range("S6").Value = 0.2390917341253824
$ws.Range("T6").Value = 0.2390917341253824

# Row 7 (FAPs -> MuSCs)
$ws.Range("I7").Value = 0.9992428949822291
$ws.Range("J7").Value = 0.9992428949822291
$ws.Range("O7").Value = 0.136569704933115
$ws.Range("P7").Value = 0.136569704933115
$ws.Range("S7").Value = 0.1364663073242347
$ws.Range("T7").Value = 0.1364663073242346

# Row 8 (MuSCs -> ECs)
$ws.Range("I8").Value = 0.0006104881997874136
$ws.Range("J8").Value = 0.0006104881997874135
$ws.Range("M8").Value = 123.2806423333333
$ws.Range("N8").Value = 369.841927
$ws.Range("O8").Value = 0.6241574062367528
$ws.Range("P8").Value = 0.6241574062367526
$ws.Range("Q8").Value = 21.32845518171044
$ws.Range("R8").Value = 191.956096635394
$ws.Range("S8").Value = 0.0003810407313174565
$ws.Range("T8").Value = 0.0003810407313174564

# Row 9 (MuSCs -> FAPs)
$ws.Range("I9").Value = 0.0006104881997874136
$ws.Range("J9").Value = 0.0006104881997874135
$ws.Range("O9").Value = 0.2392728888301323
$ws.Range("P9").Value = 0.2392728888301322
$ws.Range("S9").Value = 0.0001460732751598414
$ws.Range("T9").Value = 0.0001460732751598413

# Row 10 (MuSCs -> MuSCs)
$ws.Range("I10").Value = 0.0006104881997874136
$ws.Range("J10").Value = 0.0006104881997874135
$ws.Range("O10").Value = 0.136569704933115
$ws.Range("P10").Value = 0.136569704933115
$ws.Range("S10").Value = 0.00008337419331011565
$ws.Range("T10").Value = 0.00008337419331011563
